$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Progress" column (A) for rows 5-8 from "In Progress" to "Complete"
$ws.Range("A5").Value = "Complete"
$ws.Range("A6").Value = "Complete"
$ws.Range("A7").Value = "Complete"
$ws.Range("A8").Value = "Complete"

# Set explicit column widths for A, C, D (Employee Wireframe layout update)
# (values chosen so the engine's internal pixel-rounded ColumnWidth lands as
# close as possible to the target stored widths of 19 / 15.140625 / 65.42578125)
$ws.Range("A:A").ColumnWidth = 18.1666666666667
$ws.Range("C:C").ColumnWidth = 14.3333333333333
$ws.Range("D:D").ColumnWidth = 64.6666666666667

# Update the active selection to match the new cursor position
$ws.Range("B13").Select()
